$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from row 15's A cell (style s="1") onto A16
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1.004263002895115
$ws.Range("D16").Value = 1.015243375825277
$ws.Range("E16").Value = 0.998206045989162
$ws.Range("F16").Value = 1.004263002895115
$ws.Range("G16").Value = 1.006928091757363
$ws.Range("H16").Value = 0.9995347420427434
$ws.Range("I16").Value = 1.001212588345775
$ws.Range("J16").Value = 1.015243375825277
$ws.Range("K16").Value = 1.006724710907219
$ws.Range("L16").Value = 1.005493856901167
$ws.Range("M16").Value = 1.004231307809239
